# Set-PriceRow fixed; was using FindElementBy instead of FindElementsBy.
# Now testing Change operation including Find-Product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now exercises the "Change" operation (was "Add") for the
# Find-Product / Set-PriceRow test pass.
$ws.Range("A2").Value = "Change"

# Set-PriceRow only touches the Setup Price going forward; clear the old
# Regular Price test value and set the new Setup Price.
$ws.Range("AN2").ClearContents()
$ws.Range("AO2").Value = 0.75

# Leave the cursor on the row-3 Operation cell, matching the new test run.
$ws.Range("A3").Select()
